$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.939086437225342
$ws.Range("B1").Value = 6.344685554504395
$ws.Range("C1").Value = 3.372495174407959
$ws.Range("D1").Value = 1.488076210021973
$ws.Range("E1").Value = 1.045774817466736
